$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying source data was re-sorted by party name (case-sensitive,
# ordinal ASCII order) within each country group before the workbook was
# regenerated. For Belgium, Czech Republic and France this reshuffles which
# party's values land on which pre-existing row. Row numbers/count do not
# change -- only the party label (col B) and the 23 numeric expert-survey
# columns (C:Y) for the affected rows are updated in place to match.

# --- Belgium (BE) rows 16-18: Parti Populaire, PS, PVDA+ -> sorted ---
$ws.Range("B16").Value = 'PS'
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 2.5999999
$ws.Range("E16").Value = 2.4000001
$ws.Range("F16").Value = 3.4000001
$ws.Range("G16").Value = 2.4000001
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 1.8
$ws.Range("J16").Value = 1.75
$ws.Range("K16").Value = 2.2
$ws.Range("L16").Value = 1.6
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 2.8
$ws.Range("O16").Value = 2
$ws.Range("P16").Value = 3.3333333
$ws.Range("Q16").Value = 3.8
$ws.Range("R16").Value = 6.5999999
$ws.Range("S16").Value = 5.8000002
$ws.Range("T16").Value = 2
$ws.Range("U16").Value = 2.25
$ws.Range("V16").Value = 8.3999996
$ws.Range("W16").Value = 3.2
$ws.Range("X16").Value = 2.5999999
$ws.Range("Y16").Value = 2.25

$ws.Range("B17").Value = 'PVDA+'
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 0.40000001
$ws.Range("E17").Value = 0.2
$ws.Range("F17").Value = 2.8
$ws.Range("G17").Value = 0.2
$ws.Range("H17").Value = 0.40000001
$ws.Range("I17").Value = 0.2
$ws.Range("J17").Value = 0.25
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 2.25
$ws.Range("M17").Value = 0.80000001
$ws.Range("N17").Value = 1.8
$ws.Range("O17").Value = 1.6
$ws.Range("P17").Value = 2.6666667
$ws.Range("Q17").Value = 4.5999999
$ws.Range("R17").Value = 9.25
$ws.Range("S17").Value = 7.8000002
$ws.Range("T17").Value = 1.6
$ws.Range("U17").Value = 1
$ws.Range("V17").Value = 8.6000004
$ws.Range("W17").Value = 4.4000001
$ws.Range("X17").Value = 8.3999996
$ws.Range("Y17").Value = 6

$ws.Range("B18").Value = 'Parti Populaire'
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 7.75
$ws.Range("E18").Value = 8.5
$ws.Range("F18").Value = 7.5
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 7
$ws.Range("I18").Value = 6.5
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = 8.5
$ws.Range("L18").Value = 7
$ws.Range("M18").Value = 2
$ws.Range("N18").Value = 8.5
$ws.Range("O18").Value = 9
$ws.Range("P18").Value = 6
$ws.Range("Q18").Value = 7.5
$ws.Range("R18").Value = 8.5
$ws.Range("S18").Value = 6
$ws.Range("T18").Value = 9
$ws.Range("U18").Value = 10
$ws.Range("V18").Value = 8.8000002
$ws.Range("W18").Value = 6.1999998
$ws.Range("X18").Value = 6.5
$ws.Range("Y18").Value = 5.3333335

# --- Czech Republic (CZ) rows 33-38: re-sorted party order ---
$ws.Range("B33").Value = 'KDU-ČSL'
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 5.9285712
$ws.Range("E33").Value = 5.5714288
$ws.Range("F33").Value = 7.6428571
$ws.Range("G33").Value = 4.9000001
$ws.Range("H33").Value = 4.8333335
$ws.Range("I33").Value = 4.9166665
$ws.Range("J33").Value = 4.7692308
$ws.Range("K33").Value = 7.3333335
$ws.Range("L33").Value = 8.1538458
$ws.Range("M33").Value = 9.1538458
$ws.Range("N33").Value = 7
$ws.Range("O33").Value = 7.7272725
$ws.Range("P33").Value = 8
$ws.Range("Q33").Value = 4.5833335
$ws.Range("R33").Value = 3.8181818
$ws.Range("S33").Value = 3
$ws.Range("T33").Value = 5.75
$ws.Range("U33").Value = 5.0769229
$ws.Range("V33").Value = 6.2142859
$ws.Range("W33").Value = 7.2142859
$ws.Range("X33").Value = 2.4615386
$ws.Range("Y33").Value = 5.5384617

$ws.Range("B34").Value = 'KSČM'
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 1.0714285
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = 6.5714288
$ws.Range("G34").Value = 1.0909091
$ws.Range("H34").Value = 2.0833333
$ws.Range("I34").Value = 0.84615386
$ws.Range("J34").Value = 0.61538464
$ws.Range("K34").Value = 7
$ws.Range("L34").Value = 5.6923075
$ws.Range("M34").Value = 0.46153846
$ws.Range("N34").Value = 6.6666665
$ws.Range("O34").Value = 6.8000002
$ws.Range("P34").Value = 6.0833335
$ws.Range("Q34").Value = 6.9166665
$ws.Range("R34").Value = 6.3000002
$ws.Range("S34").Value = 8.416667
$ws.Range("T34").Value = 6.181818
$ws.Range("U34").Value = 7.6153846
$ws.Range("V34").Value = 7.2857141
$ws.Range("W34").Value = 3.6153846
$ws.Range("X34").Value = 5.6923075
$ws.Range("Y34").Value = 5.4615383

$ws.Range("B35").Value = 'ODS'
$ws.Range("C35").Value = 5
$ws.Range("D35").Value = 8
$ws.Range("E35").Value = 8.1428576
$ws.Range("F35").Value = 6
$ws.Range("G35").Value = 8.181818
$ws.Range("H35").Value = 7.5
$ws.Range("I35").Value = 7.9166665
$ws.Range("J35").Value = 7.9166665
$ws.Range("K35").Value = 5.0833335
$ws.Range("L35").Value = 5.5
$ws.Range("M35").Value = 5.2307692
$ws.Range("N35").Value = 7.875
$ws.Range("O35").Value = 7.4166665
$ws.Range("P35").Value = 3.4166667
$ws.Range("Q35").Value = 7.5384617
$ws.Range("R35").Value = 6.3636365
$ws.Range("S35").Value = 2.4166667
$ws.Range("T35").Value = 6.6666665
$ws.Range("U35").Value = 7.4615383
$ws.Range("V35").Value = 7.8571429
$ws.Range("W35").Value = 5.2142859
$ws.Range("X35").Value = 2.1538463
$ws.Range("Y35").Value = 3.1538463

$ws.Range("B36").Value = 'TOP 09'
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = 7.2142859
$ws.Range("E36").Value = 7.8571429
$ws.Range("F36").Value = 5.5
$ws.Range("G36").Value = 7.909091
$ws.Range("H36").Value = 7.25
$ws.Range("I36").Value = 7.6153846
$ws.Range("J36").Value = 7.6923075
$ws.Range("K36").Value = 4.0833335
$ws.Range("L36").Value = 4.3846154
$ws.Range("M36").Value = 6.7692308
$ws.Range("N36").Value = 5
$ws.Range("O36").Value = 6
$ws.Range("P36").Value = 3.3333333
$ws.Range("Q36").Value = 6.25
$ws.Range("R36").Value = 4
$ws.Range("S36").Value = 2.25
$ws.Range("T36").Value = 4.75
$ws.Range("U36").Value = 3.7692308
$ws.Range("V36").Value = 7.7142859
$ws.Range("W36").Value = 4.5714288
$ws.Range("X36").Value = 1.9230769
$ws.Range("Y36").Value = 5.0769229

$ws.Range("B37").Value = 'Úsvit přímé demokracie Tomia Okamury'
$ws.Range("C37").Value = 7
$ws.Range("D37").Value = 7.6923075
$ws.Range("E37").Value = 5.3333335
$ws.Range("F37").Value = 7.7142859
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 4.4285712
$ws.Range("I37").Value = 4
$ws.Range("J37").Value = 5.5
$ws.Range("K37").Value = 8.583333
$ws.Range("L37").Value = 8.181818
$ws.Range("M37").Value = 3.1666667
$ws.Range("N37").Value = 9.3999996
$ws.Range("O37").Value = 9.666667
$ws.Range("P37").Value = 6.4285712
$ws.Range("Q37").Value = 6.5555553
$ws.Range("R37").Value = 4.1666665
$ws.Range("S37").Value = 6.6999998
$ws.Range("T37").Value = 9.6153851
$ws.Range("U37").Value = 9.2307692
$ws.Range("V37").Value = 4
$ws.Range("W37").Value = 6.2307692
$ws.Range("X37").Value = 9.4615383
$ws.Range("Y37").Value = 8.6153851

$ws.Range("B38").Value = 'ČSSD'
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 3.1428571
$ws.Range("E38").Value = 2.7142856
$ws.Range("F38").Value = 4.4285712
$ws.Range("G38").Value = 2.6363637
$ws.Range("H38").Value = 3.8333333
$ws.Range("I38").Value = 2.8461537
$ws.Range("J38").Value = 2.8461537
$ws.Range("K38").Value = 4.8333335
$ws.Range("L38").Value = 3.3076923
$ws.Range("M38").Value = 2.7692308
$ws.Range("N38").Value = 4.3333335
$ws.Range("O38").Value = 4.25
$ws.Range("P38").Value = 4.9166665
$ws.Range("Q38").Value = 5.6923075
$ws.Range("R38").Value = 4.5
$ws.Range("S38").Value = 4.5833335
$ws.Range("T38").Value = 4.6923075
$ws.Range("U38").Value = 5.0769229
$ws.Range("V38").Value = 7.8571429
$ws.Range("W38").Value = 3.1428571
$ws.Range("X38").Value = 1.5
$ws.Range("Y38").Value = 5.2307692

# --- France (FR) rows 83-85: re-sorted party order ---
$ws.Range("B83").Value = 'PR (Parti Radical Valoisien)'
$ws.Range("C83").Value = 3
$ws.Range("D83").Value = 6.1428571
$ws.Range("E83").Value = 6.5
$ws.Range("F83").Value = 5.8571429
$ws.Range("G83").Value = 6.3333335
$ws.Range("H83").Value = 5.8333335
$ws.Range("I83").Value = 6
$ws.Range("J83").Value = 5.8333335
$ws.Range("K83").Value = 5.6666665
$ws.Range("L83").Value = 4.8333335
$ws.Range("M83").Value = 3.8333333
$ws.Range("N83").Value = 6.4000001
$ws.Range("O83").Value = 6.1666665
$ws.Range("P83").Value = 5.4000001
$ws.Range("Q83").Value = 6.6666665
$ws.Range("R83").Value = 4.5
$ws.Range("S83").Value = 4.4000001
$ws.Range("T83").Value = 5.8000002
$ws.Range("U83").Value = 5.8333335
$ws.Range("V83").Value = 7.181818
$ws.Range("W83").Value = 5.5
$ws.Range("X83").Value = 3.5
$ws.Range("Y83").Value = 4

$ws.Range("B84").Value = 'PS (Parti Socialiste)'
$ws.Range("C84").Value = 9
$ws.Range("D84").Value = 3.8333333
$ws.Range("E84").Value = 3.8333333
$ws.Range("F84").Value = 3.3636363
$ws.Range("G84").Value = 3.8181818
$ws.Range("H84").Value = 4.2727275
$ws.Range("I84").Value = 3.090909
$ws.Range("J84").Value = 3.8181818
$ws.Range("K84").Value = 3.909091
$ws.Range("L84").Value = 1.9090909
$ws.Range("M84").Value = 1.9090909
$ws.Range("N84").Value = 4.6999998
$ws.Range("O84").Value = 4.7272725
$ws.Range("P84").Value = 3.3333333
$ws.Range("Q84").Value = 5.3636365
$ws.Range("R84").Value = 3.3
$ws.Range("S84").Value = 3.3
$ws.Range("T84").Value = 3.7272727
$ws.Range("U84").Value = 4.2727275
$ws.Range("V84").Value = 7.0833335
$ws.Range("W84").Value = 6.5454545
$ws.Range("X84").Value = 3.2
$ws.Range("Y84").Value = 3.7

$ws.Range("B85").Value = 'Parti Radical de Gauche'
$ws.Range("C85").Value = 7
$ws.Range("D85").Value = 3.8181818
$ws.Range("E85").Value = 3.75
$ws.Range("F85").Value = 3.4166667
$ws.Range("G85").Value = 3.4444444
$ws.Range("H85").Value = 3.8
$ws.Range("I85").Value = 2.7777777
$ws.Range("J85").Value = 3.5454545
$ws.Range("K85").Value = 3.625
$ws.Range("L85").Value = 2
$ws.Range("M85").Value = 1.7272727
$ws.Range("N85").Value = 4
$ws.Range("O85").Value = 5
$ws.Range("P85").Value = 5.1666665
$ws.Range("Q85").Value = 5.3000002
$ws.Range("R85").Value = 3.3333333
$ws.Range("S85").Value = 4.25
$ws.Range("T85").Value = 4
$ws.Range("U85").Value = 4
$ws.Range("V85").Value = 6.6363635
$ws.Range("W85").Value = 6.4545455
$ws.Range("X85").Value = 3.7777777
$ws.Range("Y85").Value = 3.5
